$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 3.075165666666667
$ws.Range("H2").Value = 9.225497000000001
$ws.Range("I2").Value = 0.02641273658732285
$ws.Range("J2").Value = 0.02641273658732285
$ws.Range("M2").Value = 721.1356606666667
$ws.Range("N2").Value = 2163.406982
$ws.Range("O2").Value = 0.8508208584329936
$ws.Range("P2").Value = 0.8508208584329936
$ws.Range("Q2").Value = 2217.611624691117
$ws.Range("R2").Value = 19958.50462222005
$ws.Range("S2").Value = 0.02247250721679056
$ws.Range("T2").Value = 0.02247250721679056
$ws.Range("G3").Value = 3.075165666666667
$ws.Range("H3").Value = 9.225497000000001
$ws.Range("I3").Value = 0.02641273658732285
$ws.Range("J3").Value = 0.02641273658732285
$ws.Range("O3").Value = 0.002793596814304166
$ws.Range("P3").Value = 0.002793596814304166
$ws.Range("Q3").Value = 7.281336263324445
$ws.Range("R3").Value = 65.53202636992
$ws.Range("S3").Value = 0.00007378653678740019
$ws.Range("T3").Value = 0.00007378653678740019
$ws.Range("G4").Value = 3.075165666666667
$ws.Range("H4").Value = 9.225497000000001
$ws.Range("I4").Value = 0.02641273658732285
$ws.Range("J4").Value = 0.02641273658732285
$ws.Range("O4").Value = 0.1463855447527022
$ws.Range("P4").Value = 0.1463855447527022
$ws.Range("Q4").Value = 381.5448134736815
$ws.Range("R4").Value = 3433.903321263133
$ws.Range("S4").Value = 0.003866442833744882
$ws.Range("T4").Value = 0.003866442833744883
$ws.Range("I5").Value = 0.549422396165273
$ws.Range("J5").Value = 0.5494223961652731
$ws.Range("M5").Value = 721.1356606666667
$ws.Range("N5").Value = 2163.406982
$ws.Range("O5").Value = 0.8508208584329936
$ws.Range("P5").Value = 0.8508208584329936
$ws.Range("Q5").Value = 46129.46820461415
$ws.Range("R5").Value = 415165.2138415273
$ws.Range("S5").Value = 0.4674600347476499
$ws.Range("T5").Value = 0.46746003474765
$ws.Range("I6").Value = 0.549422396165273
$ws.Range("J6").Value = 0.5494223961652731
$ws.Range("O6").Value = 0.002793596814304166
$ws.Range("P6").Value = 0.002793596814304166
$ws.Range("S6").Value = 0.001534864655634668
$ws.Range("T6").Value = 0.001534864655634668
$ws.Range("I7").Value = 0.549422396165273
$ws.Range("J7").Value = 0.5494223961652731
$ws.Range("O7").Value = 0.1463855447527022
$ws.Range("P7").Value = 0.1463855447527022
$ws.Range("S7").Value = 0.08042749676198842
$ws.Range("T7").Value = 0.08042749676198846
$ws.Range("I8").Value = 0.424164867247404
$ws.Range("J8").Value = 0.4241648672474041
$ws.Range("M8").Value = 721.1356606666667
$ws.Range("N8").Value = 2163.406982
$ws.Range("O8").Value = 0.8508208584329936
$ws.Range("P8").Value = 0.8508208584329936
$ws.Range("Q8").Value = 35612.85432441247
$ws.Range("R8").Value = 320515.6889197122
$ws.Range("S8").Value = 0.3608883164685531
$ws.Range("T8").Value = 0.3608883164685531
$ws.Range("I9").Value = 0.424164867247404
$ws.Range("J9").Value = 0.4241648672474041
$ws.Range("O9").Value = 0.002793596814304166
$ws.Range("P9").Value = 0.002793596814304166
$ws.Range("S9").Value = 0.001184945621882097
$ws.Range("T9").Value = 0.001184945621882098
$ws.Range("I10").Value = 0.424164867247404
$ws.Range("J10").Value = 0.4241648672474041
$ws.Range("O10").Value = 0.1463855447527022
$ws.Range("P10").Value = 0.1463855447527022
$ws.Range("S10").Value = 0.06209160515696883
$ws.Range("T10").Value = 0.06209160515696886
